$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  18 = @(' Dubai (DSC)', ' October 13 2020', 'Super Kings won by 20 runs', 'Sunrisers Hyderabad', 'Chennai Super Kings', 'David Warner (c)', '9', '13', '0', '0', '69.23')
  19 = @(' Dubai (DSC)', ' October 22 2020', 'Sunrisers won by 8 wickets (with 11 balls remaining)', 'Sunrisers Hyderabad', 'Rajasthan Royals', 'David Warner (c)', '4', '4', '1', '0', '100.00')
  20 = @(' Abu Dhabi', ' October 18 2020', 'Match tied (KKR won the one-over eliminator)', 'Sunrisers Hyderabad', 'Kolkata Knight Riders', 'David Warner (c)', '47', '33', '5', '0', '142.42')
  21 = @(' Dubai (DSC)', ' October 27 2020', 'Sunrisers won by 88 runs', 'Sunrisers Hyderabad', 'Delhi Capitals', 'David Warner (c)', '66', '34', '8', '2', '194.11')
  22 = @(' Dubai (DSC)', ' October 02 2020', 'Sunrisers won by 7 runs', 'Sunrisers Hyderabad', 'Chennai Super Kings', 'David Warner (c)', '28', '29', '3', '0', '96.55')
  23 = @(' Abu Dhabi', ' September 29 2020', 'Sunrisers won by 15 runs', 'Sunrisers Hyderabad', 'Delhi Capitals', 'David Warner (c)', '45', '33', '3', '2', '136.36')
  24 = @(' Sharjah', ' October 04 2020', 'Mumbai won by 34 runs', 'Sunrisers Hyderabad', 'Mumbai Indians', 'David Warner (c)', '60', '44', '5', '2', '136.36')
  25 = @(' Sharjah', ' November 03 2020', 'Sunrisers won by 10 wickets (with 17 balls remaining)', 'Sunrisers Hyderabad', 'Mumbai Indians', 'David Warner (c)', '85', '58', '10', '1', '146.55')
  26 = @(' Abu Dhabi', ' September 26 2020', 'KKR won by 7 wickets (with 12 balls remaining)', 'Sunrisers Hyderabad', 'Kolkata Knight Riders', 'David Warner (c)', '36', '30', '2', '1', '120.00')
  27 = @(' Abu Dhabi', ' November 06 2020', 'Sunrisers won by 6 wickets (with 2 balls remaining)', 'Sunrisers Hyderabad', 'Royal Challengers Bangalore', 'David Warner (c)', '17', '17', '3', '0', '100.00')
  28 = @(' Sharjah', ' October 31 2020', 'Sunrisers won by 5 wickets (with 35 balls remaining)', 'Sunrisers Hyderabad', 'Royal Challengers Bangalore', 'David Warner (c)', '8', '5', '0', '1', '160.00')
  29 = @(' Dubai (DSC)', ' September 21 2020', 'RCB won by 10 runs', 'Sunrisers Hyderabad', 'Royal Challengers Bangalore', 'David Warner (c)', '6', '6', '1', '0', '100.00')
  30 = @(' Abu Dhabi', ' November 08 2020', 'Capitals won by 17 runs', 'Sunrisers Hyderabad', 'Delhi Capitals', 'David Warner (c)', '2', '3', '0', '0', '66.66')
  31 = @(' Dubai (DSC)', ' October 11 2020', 'Royals won by 5 wickets (with 1 ball remaining)', 'Sunrisers Hyderabad', 'Rajasthan Royals', 'David Warner (c)', '48', '38', '3', '2', '126.31')
  32 = @(' Dubai (DSC)', ' October 24 2020', 'Kings XI won by 12 runs', 'Sunrisers Hyderabad', 'Kings XI Punjab', 'David Warner (c)', '35', '20', '3', '2', '175.00')
  33 = @(' Dubai (DSC)', ' October 08 2020', 'Sunrisers won by 69 runs', 'Sunrisers Hyderabad', 'Kings XI Punjab', 'David Warner (c)', '52', '40', '5', '1', '130.00')
}

$cols = @("A","B","C","D","E","F","G","H","I","J","K")
foreach ($r in $data.Keys) {
  $vals = $data[$r]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $cell = $ws.Range($cols[$i] + $r)
    $cell.NumberFormat = "@"
    $cell.Value = $vals[$i]
  }
}
